$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 89; this shifts the existing rows 89..188
# (and their formatting, e.g. the date style on column D) down to 90..189,
# and grows the sheet dimension to A1:R189 automatically.
$ws.Rows.Item(89).Insert()

# Populate the newly inserted row 89 with the new weekly price record.
$ws.Cells.Item(89, 1).Value = 10
$ws.Cells.Item(89, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(89, 3).Value = "La Araucanía"
$ws.Cells.Item(89, 4).Value = 44810
$ws.Cells.Item(89, 5).Value = 9
$ws.Cells.Item(89, 6).Value = 100112012
$ws.Cells.Item(89, 7).Value = "Espinaca"
$ws.Cells.Item(89, 8).Value = "Sin especificar"
$ws.Cells.Item(89, 9).Value = "Primera"
$ws.Cells.Item(89, 10).Value = 40
$ws.Cells.Item(89, 11).Value = 10000
$ws.Cells.Item(89, 12).Value = 10000
$ws.Cells.Item(89, 13).Value = 10000
$ws.Cells.Item(89, 14).Value = "`$/docena de atados"
$ws.Cells.Item(89, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(89, 16).Value = 3333
$ws.Cells.Item(89, 17).Value = 3
$ws.Cells.Item(89, 18).Value = "Hortaliza"
